# Rename the two original sheets.
$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Description").Name   = "experiment-description"
$wb.Worksheets.Item("Specification").Name = "experiment-specification"

# Add the two new sheets. Worksheets.Add() inserts right before the
# currently-active sheet ("experiment-specification"), so adding
# "run-specification" first and "run-description" second makes
# "run-specification" get the lower internal sheetId (created first)
# while ending up to the right of "run-description" once both are
# relocated to the end of the tab strip below.
$wsRunSpec = $wb.Worksheets.Add()
$wsRunSpec.Name = "run-specification"

$wsRunDesc = $wb.Worksheets.Add()
$wsRunDesc.Name = "run-description"

# Move both new sheets to the very end, "run-description" first so it
# lands before "run-specification" in the final tab order:
# experiment-description, experiment-specification, run-description, run-specification
$wsRunDesc.Move([Type]::Missing, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsRunSpec.Move([Type]::Missing, $wb.Worksheets.Item($wb.Worksheets.Count))

# Re-resolve fresh references by name: the object handles obtained before
# Move() no longer point at the right sheet once the tab order changes.
$wsRunDesc = $wb.Worksheets.Item("run-description")
$wsRunSpec = $wb.Worksheets.Item("run-specification")

# Populate "run-description": run_id / name / description.
$wsRunDesc.Range("A1").Value = "run_id"
$wsRunDesc.Range("B1").Value = "name"
$wsRunDesc.Range("C1").Value = "description"
$wsRunDesc.Range("A2").Value = 1
$wsRunDesc.Range("B2").Value = "Base"
$wsRunDesc.Range("C2").Value = "Dummy"
$wsRunDesc.Range("E17").Select()

# Populate "run-specification": run_id / experiment_id.
$wsRunSpec.Range("A1").Value = "run_id"
$wsRunSpec.Range("B1").Value = "experiment_id"
$wsRunSpec.Range("A2").Value = 1
$wsRunSpec.Range("B2").Value = 1
$wsRunSpec.Columns.Item(2).ColumnWidth = 14.5
$wsRunSpec.Range("D18").Select()

# experiment-description: widen column A, select B1:C2.
$wsExpDesc = $wb.Worksheets.Item("experiment-description")
$wsExpDesc.Columns.Item(1).ColumnWidth = 14
$wsExpDesc.Range("A1").Value = "experiment_id"
$wsExpDesc.Range("B1:C2").Select()

# experiment-specification: rename id column header, move selection to A2
# (this sheet is no longer the active/tabSelected one).
$wsExpSpec = $wb.Worksheets.Item("experiment-specification")
$wsExpSpec.Range("A1").Value = "experiment_id"
$wsExpSpec.Range("A2").Select()

# Make "run-specification" the active sheet/tab, matching activeTab="3".
$wsRunSpec.Activate()
